$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026474207863608
$ws.Range("D2").Value = 1.031864070389793
$ws.Range("E2").Value = 1.049066468113903
$ws.Range("F2").Value = 1.053548558217059
$ws.Range("I2").Value = 1.03432358770371
$ws.Range("J2").Value = 1.031637448183044
$ws.Range("K2").Value = 1.034670957818054
$ws.Range("L2").Value = 1.051824561480208
$ws.Range("M2").Value = 1.056294221939822
$ws.Range("N2").Value = 1.014595877209362

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027247949604007
$ws.Range("D3").Value = 1.032445002752774
$ws.Range("E3").Value = 1.050259521086883
$ws.Range("F3").Value = 1.05479477476538
$ws.Range("I3").Value = 1.034494033626405
$ws.Range("J3").Value = 1.032051870079746
$ws.Range("K3").Value = 1.035061152948912
$ws.Range("L3").Value = 1.052828801864592
$ws.Range("M3").Value = 1.057352396207918
$ws.Range("N3").Value = 1.014732421722262

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027748699526991
$ws.Range("D4").Value = 1.03282082991619
$ws.Range("E4").Value = 1.05103244487655
$ws.Range("F4").Value = 1.055602063799782
$ws.Range("I4").Value = 1.034602914682867
$ws.Range("J4").Value = 1.032319442113542
$ws.Range("K4").Value = 1.035312861085237
$ws.Range("L4").Value = 1.053478979173613
$ws.Range("M4").Value = 1.058037447087229
$ws.Range("N4").Value = 1.014820577475987

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027959234343575
$ws.Range("D5").Value = 1.032978808199847
$ws.Range("E5").Value = 1.051357607312648
$ws.Range("F5").Value = 1.055941664940318
$ws.Range("I5").Value = 1.034648350334723
$ws.Range("J5").Value = 1.032431788045931
$ws.Range("K5").Value = 1.035418492810136
$ws.Range("L5").Value = 1.053752401911521
$ws.Range("M5").Value = 1.058325524538406
$ws.Range("N5").Value = 1.014857590514149

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027994585172886
$ws.Range("D6").Value = 1.033005332263618
$ws.Range("E6").Value = 1.051412216742523
$ws.Range("F6").Value = 1.055998698196115
$ws.Range("I6").Value = 1.034655959344111
$ws.Range("J6").Value = 1.032450643107054
$ws.Range("K6").Value = 1.035436217882403
$ws.Range("L6").Value = 1.053798315976784
$ws.Range("M6").Value = 1.058373898832133
$ws.Range("N6").Value = 1.014863802363227

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027751512628389
$ws.Range("D7").Value = 1.032822940908116
$ws.Range("E7").Value = 1.051036788827068
$ws.Range("F7").Value = 1.055606600714195
$ws.Range("I7").Value = 1.034603523124991
$ws.Range("J7").Value = 1.032320943841941
$ws.Range("K7").Value = 1.035314273275141
$ws.Range("L7").Value = 1.053482632314519
$ws.Range("M7").Value = 1.058041296069351
$ws.Range("N7").Value = 1.014821072233609

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026735678627482
$ws.Range("D8").Value = 1.032060413802893
$ws.Range("E8").Value = 1.049469471902185
$ws.Range("F8").Value = 1.053969536142555
$ws.Range("I8").Value = 1.034381482016016
$ws.Range("J8").Value = 1.031777624927771
$ws.Range("K8").Value = 1.034802985775813
$ws.Range("L8").Value = 1.052163873561381
$ws.Range("M8").Value = 1.056651766406289
$ws.Range("N8").Value = 1.014642063807047

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024946377031185
$ws.Range("D9").Value = 1.030716236682505
$ws.Range("E9").Value = 1.046714824649397
$ws.Range("F9").Value = 1.051091716817481
$ws.Range("I9").Value = 1.033979456214411
$ws.Range("J9").Value = 1.030815774451252
$ws.Range("K9").Value = 1.033896147666193
$ws.Range("L9").Value = 1.049842840946168
$ws.Range("M9").Value = 1.05420583685175
$ws.Range("N9").Value = 1.014325128039079

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023754077491781
$ws.Range("D10").Value = 1.029819867400887
$ws.Range("E10").Value = 1.044883174593509
$ws.Range("F10").Value = 1.049177773666509
$ws.Range("I10").Value = 1.033704242765432
$ws.Range("J10").Value = 1.03017159984419
$ws.Range("K10").Value = 1.033287693160913
$ws.Range("L10").Value = 1.048297332003471
$ws.Range("M10").Value = 1.052576941008561
$ws.Range("N10").Value = 1.014112847526259

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023237949731908
$ws.Range("D11").Value = 1.029431688769382
$ws.Range("E11").Value = 1.044091175195591
$ws.Range("F11").Value = 1.048350098447602
$ws.Range("I11").Value = 1.033583373485185
$ws.Range("J11").Value = 1.029891978545916
$ws.Range("K11").Value = 1.033023314171897
$ws.Range("L11").Value = 1.047628540328871
$ws.Range("M11").Value = 1.051872013945957
$ws.Range("N11").Value = 1.014020696746699

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023046260086721
$ws.Range("D12").Value = 1.029287496702693
$ws.Range("E12").Value = 1.043797158350464
$ws.Range("F12").Value = 1.048042823667393
$ws.Range("I12").Value = 1.033538222454792
$ws.Range("J12").Value = 1.029788012017042
$ws.Range("K12").Value = 1.032924975681148
$ws.Range("L12").Value = 1.047380184309569
$ws.Range("M12").Value = 1.051610231463272
$ws.Range("N12").Value = 1.013986433300672

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02308737709661
$ws.Range("D13").Value = 1.029318426603637
$ws.Range("E13").Value = 1.043860218418452
$ws.Range("F13").Value = 1.048108727894364
$ws.Range("I13").Value = 1.033547919024382
$ws.Range("J13").Value = 1.029810317825432
$ws.Range("K13").Value = 1.032946075776212
$ws.Range("L13").Value = 1.047433454679661
$ws.Range("M13").Value = 1.051666382042786
$ws.Range("N13").Value = 1.013993784484524

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023222104125671
$ws.Range("D14").Value = 1.029419769906431
$ws.Range("E14").Value = 1.044066868269675
$ws.Range("F14").Value = 1.048324695746
$ws.Range("I14").Value = 1.033579646481183
$ws.Range("J14").Value = 1.029883386730917
$ws.Range("K14").Value = 1.033015188257669
$ws.Range("L14").Value = 1.047608009861954
$ws.Range("M14").Value = 1.051850373728009
$ws.Range("N14").Value = 1.01401786522251

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023305117005359
$ws.Range("D15").Value = 1.029482210210357
$ws.Range("E15").Value = 1.044194214166587
$ws.Range("F15").Value = 1.048457781930056
$ws.Range("I15").Value = 1.033599161069594
$ws.Range("J15").Value = 1.02992839330809
$ws.Range("K15").Value = 1.033057752698406
$ws.Range("L15").Value = 1.047715567396345
$ws.Range("M15").Value = 1.051963744853002
$ws.Range("N15").Value = 1.014032697591539

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023788334372642
$ws.Range("D16").Value = 1.029845628721646
$ws.Range("E16").Value = 1.04493576045014
$ws.Range("F16").Value = 1.049232726318828
$ws.Range("I16").Value = 1.033712228696011
$ws.Range("J16").Value = 1.0301901429359
$ws.Range("K16").Value = 1.033305219946978
$ws.Range("L16").Value = 1.048341726375507
$ws.Range("M16").Value = 1.052623732977826
$ws.Range("N16").Value = 1.014118958409367

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024091483801592
$ws.Range("D17").Value = 1.030073580610195
$ws.Range("E17").Value = 1.045401211339791
$ws.Range("F17").Value = 1.04971911562757
$ws.Range("I17").Value = 1.033782698315424
$ws.Range("J17").Value = 1.030354147581318
$ws.Range("K17").Value = 1.033460205557853
$ws.Range("L17").Value = 1.048734612354271
$ws.Range("M17").Value = 1.053037831564425
$ws.Range("N17").Value = 1.014173005662409

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024268319717368
$ws.Range("D18").Value = 1.030206536666077
$ws.Range("E18").Value = 1.045672808815397
$ws.Range("F18").Value = 1.050002922039653
$ws.Range("I18").Value = 1.033823637981795
$ws.Range("J18").Value = 1.030449742138028
$ws.Range("K18").Value = 1.033550517772464
$ws.Range("L18").Value = 1.048963817044998
$ws.Range("M18").Value = 1.053279406473439
$ws.Range("N18").Value = 1.014204508114086

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024328618527561
$ws.Range("D19").Value = 1.030251870449292
$ws.Range("E19").Value = 1.045765434854332
$ws.Range("F19").Value = 1.05009971031627
$ws.Range("I19").Value = 1.033837569502175
$ws.Range("J19").Value = 1.030482326092689
$ws.Range("K19").Value = 1.033581296901725
$ws.Range("L19").Value = 1.04904197688793
$ws.Range("M19").Value = 1.053361783774402
$ws.Range("N19").Value = 1.014215245828034

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02405895728127
$ws.Range("D20").Value = 1.030049123976659
$ws.Range("E20").Value = 1.045351261715206
$ws.Range("F20").Value = 1.049666919931592
$ws.Range("I20").Value = 1.033775154558299
$ws.Range("J20").Value = 1.030336558315304
$ws.Range("K20").Value = 1.03344358619165
$ws.Range("L20").Value = 1.048692455176619
$ws.Range("M20").Value = 1.052993398757209
$ws.Range("N20").Value = 1.014167209218071

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0231824297233
$ws.Range("D21").Value = 1.029389926962337
$ws.Range("E21").Value = 1.044006010435985
$ws.Range("F21").Value = 1.048261094168483
$ws.Range("I21").Value = 1.033570310564689
$ws.Range("J21").Value = 1.029861872588424
$ws.Range("K21").Value = 1.032994840110878
$ws.Range("L21").Value = 1.047556605978288
$ws.Range("M21").Value = 1.051796191143921
$ws.Range("N21").Value = 1.01401077499654

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022631456833998
$ws.Range("D22").Value = 1.028975434058025
$ws.Range("E22").Value = 1.043161163331399
$ws.Range("F22").Value = 1.047378124770875
$ws.Range("I22").Value = 1.033440043136133
$ws.Range("J22").Value = 1.029562824896136
$ws.Range("K22").Value = 1.032711907169759
$ws.Range("L22").Value = 1.046842816584877
$ws.Range("M22").Value = 1.051043799186007
$ws.Range("N22").Value = 1.013912218882119

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022923524785498
$ws.Range("D23").Value = 1.029195166895507
$ws.Range("E23").Value = 1.04360894150395
$ws.Range("F23").Value = 1.047846115765801
$ws.Range("I23").Value = 1.033509239866799
$ws.Range("J23").Value = 1.029721411750252
$ws.Range("K23").Value = 1.032861969686099
$ws.Range("L23").Value = 1.047221175569109
$ws.Range("M23").Value = 1.051442624376573
$ws.Range("N23").Value = 1.01396448417367

$ws.Range("B24").Value = 1.019999999999999
$ws.Range("C24").Value = 1.02407365457263
$ws.Range("D24").Value = 1.030060174890467
$ws.Range("E24").Value = 1.045373831466977
$ws.Range("F24").Value = 1.049690504598972
$ws.Range("I24").Value = 1.03377856376443
$ws.Range("J24").Value = 1.030344506353285
$ws.Range("K24").Value = 1.033451096040699
$ws.Range("L24").Value = 1.04871150406379
$ws.Range("M24").Value = 1.053013475911887
$ws.Range("N24").Value = 1.014169828450991

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025408860006004
$ws.Range("D25").Value = 1.031063789503175
$ws.Range("E25").Value = 1.047426121602261
$ws.Range("F25").Value = 1.051834888267205
$ws.Range("I25").Value = 1.034084660747488
$ws.Range("J25").Value = 1.031064958094759
$ws.Range("K25").Value = 1.034131278122066
$ws.Range("L25").Value = 1.050442556021786
$ws.Range("M25").Value = 1.054837862709662
$ws.Range("N25").Value = 1.01440723938043
